$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new effector row (row 6) - most effectors randomly generated
$ws.Range("A6").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 500

# Move the active selection to E4
$ws.Range("E4").Select()
